# Update crypto price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.734.96'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '1.863.47'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('D4').Value = "'1.036"
$ws.Range('E4').Value = '  +1.07%  '
$ws.Range('D5').Value = "'323.15"
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').Value = "'1.031"
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('D7').Value = "'0.4422"
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('D8').Value = "'0.3800"
$ws.Range('E8').Value = '  +2.12%  '
$ws.Range('D9').Value = "'0.07462"
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('D11').Value = "'21.74"
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').Value = '1.876.51'
$ws.Range('E12').Value = '  -7.29%  '
$ws.Range('D13').Value = "'5.550"
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').Value = "'6.772"
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('D15').Value = "'0.07202"
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').Value = "'84.24"
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('D17').Value = "'1.037"
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('D18').Value = "'0.000009122"
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = "'1.031"
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').Value = "'15.55"
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = '27.741.49'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').Value = "'5.309"
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').Value = '2.091.79'
$ws.Range('E24').Value = '  -5.88%  '
$ws.Range('D25').Value = "'2.018"
$ws.Range('E25').Value = '  +5.88%  '
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').Value = "'18.85"
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('E28').Value = '  +3.60%  '
$ws.Range('D29').Value = "'5.338"
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').Value = "'118.08"
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('D31').Value = "'0.09039"
$ws.Range('D32').Value = "'1.229"
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').Value = "'0.7801"
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('D34').Value = "'3.026"
$ws.Range('E34').Value = '  +6.37%  '
$ws.Range('D35').Value = "'4.581"
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('D36').Value = "'1.033"
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').Value = "'1.151"
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').Value = "'2.881"
$ws.Range('E40').Value = '  +2.57%  '
$ws.Range('D41').Value = "'0.5227"
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('D42').Value = "'0.1695"
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('D43').Value = "'6.889"
$ws.Range('E43').Value = '  +5.53%  '
$ws.Range('D44').Value = "'8.689"
$ws.Range('E44').Value = '  +2.63%  '
$ws.Range('D45').Value = "'110.77"
$ws.Range('E45').Value = '  +2.48%  '
$ws.Range('D46').Value = "'0.06749"
$ws.Range('E46').Value = '  +7.43%  '
$ws.Range('D47').Value = "'10.66"
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('D49').Value = "'0.4726"
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = "'1.918"
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('D51').Value = "'39.78"
$ws.Range('E51').Value = '  +1.44%  '
